# Weekly update: a new price record (the most recent week) is inserted at
# row 184, pushing all the previously-existing rows (old 184..281) down by
# one (new 185..282).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 184 - this shifts rows 184-281 down to
# 185-282 and extends the used range to A1:R282.
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new weekly data point.
$ws.Cells.Item(184, 1).Value = 10
$ws.Cells.Item(184, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(184, 3).Value = "La Araucanía"
$ws.Cells.Item(184, 4).Value = 45029
$ws.Cells.Item(184, 5).Value = 9
$ws.Cells.Item(184, 6).Value = 100112005
$ws.Cells.Item(184, 7).Value = "Puerro"
$ws.Cells.Item(184, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 50
$ws.Cells.Item(184, 11).Value = 12000
$ws.Cells.Item(184, 12).Value = 12000
$ws.Cells.Item(184, 13).Value = 12000
$ws.Cells.Item(184, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(184, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(184, 16).Value = 1000
$ws.Cells.Item(184, 17).Value = 12
$ws.Cells.Item(184, 18).Value = "Hortaliza"
